# Populate the previously-empty "Web page" (column H) cells for PC members
# with their personal/profile web page URLs, and complete the associated
# sheet formatting (column widths, zoom, selection) to match the finished
# committee sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    4 = 'https://scholar.google.it/citations?user=8Z20aIMAAAAJ&hl=it'
    8 = 'https://apps.ualberta.ca/directory/person/bezemer'
    9 = 'https://www.monash.edu/it/humanise-lab/about-us/our-students/current-monash-students/tingtin-bi'
    11 = 'https://research.jetbrains.org/researchers/timofey.bryksin/'
    20 = 'https://jinfuchen.github.io/jinfu/'
    22 = 'https://cs.nju.edu.cn/chenlin/index.htm'
    24 = 'https://sites.google.com/site/ejchoi710/'
    25 = 'https://informatics.tuwien.ac.at/people/juergen-cito'
    27 = 'https://www.cs.usask.ca/faculty/zadiacodabux/'
    28 = 'https://wwwfr.uni.lu/layout/set/print/snt/people/maxime_cordy'
    30 = 'https://www.brunel.ac.uk/people/steve-counsell'
    31 = 'https://diegoeliascosta.github.io/'
    32 = 'https://diegoeliascosta.github.io/'
    39 = 'https://www.uni-trier.de/index.php?id=3569'
    41 = 'https://xiaoningdu.github.io/'
    44 = 'https://www.cristal.univ-lille.fr/~etien/'
    45 = 'https://cse.umn.edu/cs/mattia-fazzini'
    46 = 'https://fengyang-nju.github.io/'
    48 = 'https://dfucci.github.io/'
    55 = 'https://sites.google.com/view/sanondagupta147/sanonda-datta-gupta'
    56 = 'https://gyimothy.sed.hu/'
    57 = 'https://habchisarra.github.io/'
    60 = 'https://pinjiahe.github.io/'
    61 = 'https://pinjiahe.github.io/'
    62 = 'http://hematimoghadam.faculty.vru.ac.ir/?part=menu&inc=menu&id=118'
    63 = 'https://cmps.ok.ubc.ca/about/contact/fatemeh-hendijani-fard/'
    68 = 'https://xing-hu.github.io/'
    72 = 'http://sei.buaa.edu.cn/info/1082/1121.htm'
    74 = 'NULL'
    77 = 'https://kanghj.github.io/'
    80 = 'https://mkechagia.github.io/'
    83 = 'https://researchers.uq.edu.au/researcher/30167'
    84 = 'https://vovak.me/'
    85 = 'https://vovak.me/'
    87 = 'https://raux.github.io/'
    88 = 'https://lamothemax.github.io/'
    89 = 'https://lamothemax.github.io/'
    91 = 'http://www.cse.cqu.edu.cn/info/2096/4505.htm'
    93 = 'http://sei.pku.edu.cn/~lige/'
    94 = 'https://lilicoding.github.io/'
    95 = 'https://lilicoding.github.io/'
    96 = 'https://lx0704.github.io/'
    97 = 'https://binlin.info/'
    99 = 'https://liuhuigmail.github.io/'
    100 = 'https://person.zju.edu.cn/en/liuzhongxin'
    101 = 'https://person.zju.edu.cn/en/liuzhongxin'
    104 = 'https://yilinglou.github.io/'
    105 = 'https://www4.comp.polyu.edu.hk/~csxluo/'
    111 = 'https://mkaouer.net/author/mohamed-wiem-mkaouer/'
    112 = 'http://www.ranmo.site/'
    121 = 'https://sites.google.com/site/ouniaali/home'
    124 = 'https://lucapascarella.com/'
    128 = 'https://collab.di.uniba.it/luigi-quaranta/'
    136 = 'https://www.fujitsu.com/us/about/businesspolicy/tech/rd/research-staff/ripon.html'
    137 = 'https://msayagh.github.io/'
    143 = 'https://allisonius.github.io/'
    149 = 'https://scholar.google.com/citations?user=PIBtljkAAAAJ&hl=zh-CN'
    150 = 'https://home.cse.ust.hk/~shuaiw/'
    152 = 'https://wangying-neu.github.io/'
    153 = 'https://people.ucas.ac.cn/~weijun'
    156 = 'https://faculty.stevens.edu/lxiao6'
    158 = 'https://xiaoyuanxie.github.io/'
    159 = 'https://xiaoyuanxie.github.io/'
    161 = 'http://www.cse.cqu.edu.cn/info/2095/5040.htm'
    162 = 'http://www.cse.cqu.edu.cn/info/2095/5040.htm'
    163 = 'https://jinqiuyang.github.io/'
    164 = 'https://jinqiuyang.github.io/'
    165 = 'https://jinqiuyang.github.io/'
    168 = 'https://scholar.google.com/citations?user=-qRba7AAAAAJ&hl=it'
    178 = 'https://www.cs.sjtu.edu.cn/~zhonghao/'
    179 = 'http://yiikou.com/jiayuan/about/'
    180 = 'https://www.eecg.utoronto.ca/~shuruiz/'
    181 = 'http://faculty.nuaa.edu.cn/zhouyu/en/index.htm'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 8).Value = $values[$row]
}

# Widen column H (Web page) now that it holds full URLs, and nudge the
# neighbouring columns (I, J, L) to their final widths.
$ws.Columns.Item(8).ColumnWidth = 133.666666666667
$ws.Columns.Item(9).ColumnWidth = 24.8333333333333
$ws.Columns.Item(10).ColumnWidth = 17.5
$ws.Columns.Item(12).ColumnWidth = 97.6666666666667

# Reset the view: zoom back out, scroll to the top-left, and leave the
# selection on the last completed row.
$excel.ActiveWindow.Zoom = 130
$ws.Range("A1").Select()
$ws.Range("H181").Select()
